# Generate Report for Handback
#
# The "4503c4c5-86d4-4d3b-9b35-781577df6db2.md" file has now been handed
# back (it's in sync with en-US), so its status moves from
# "Ready for handoff" to "Handed back: in sync with en-US" on the Overview
# sheet (both locale columns) as well as on each locale's detail sheet.
# The "Latest Handback DateTime" for both files is refreshed with the
# timestamp of this handback run.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $handedBack
$overview.Range("C3").Value = $handedBack

# --- zh-cn detail sheet -----------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-03-14 09:52:40"
$zhcn.Range("C3").Value = $handedBack
$zhcn.Range("H3").Value = "2016-03-14 09:52:40"

# --- de-de detail sheet -----------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-03-14 09:52:54"
$dede.Range("C3").Value = $handedBack
$dede.Range("H3").Value = "2016-03-14 09:52:54"
